# YTI-1048: Bugfixes to codescheme language impl and related integration
# test content fixes.
#
# Adds a LANGUAGECODE column to the CodeSchemes sheet (between
# CLASSIFICATION and STARTDATE), populates the header + the test row value
# "fi;sv;en", and makes the CodeSchemes sheet the active/selected sheet
# (instead of Extensions_test), with the selection resting on the new
# I2 cell.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("CodeSchemes")

# Insert a new column at I, shifting the former I:L (STARTDATE, ENDDATE,
# CODESSHEET, EXTENSIONSCHEMESSHEET) right to J:M.
$ws1.Columns.Item(9).Insert()

# Populate the new LANGUAGECODE column.
$ws1.Cells.Item(1, 9).Value = "LANGUAGECODE"
$ws1.Cells.Item(2, 9).Value = "fi;sv;en"

# Match the column width used for the other "bestFit"-less text columns.
$ws1.Columns.Item(9).ColumnWidth = 13.666666666666666

# CodeSchemes becomes the active sheet (tabSelected) with I2 selected;
# this also clears tabSelected on the previously-active Extensions_test
# sheet.
[void]$ws1.Activate()
$ws1.Range("I2").Select() | Out-Null
